# Auto-generated: applies the scheduled-runner value refresh to the Leve profit tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value2 = 45.857143
$ws.Range("I11").Value2 = 45.857143
$ws.Range("K11").Value2 = 45.857143
$ws.Range("M11").Value2 = 94.14285699999999
$ws.Range("H17").Value2 = 1278.6316
$ws.Range("J17").Value2 = 1305.5428
$ws.Range("L17").Value2 = 3916.6284
$ws.Range("N17").Value2 = -4252.6284
$ws.Range("H29").Value2 = 910
$ws.Range("I29").Value2 = 325.7143
$ws.Range("J29").Value2 = 5000
$ws.Range("K29").Value2 = 977.1428999999999
$ws.Range("L29").Value2 = 15000
$ws.Range("M29").Value2 = -696.1428999999999
$ws.Range("N29").Value2 = -15562
$ws.Range("H33").Value2 = 251.66667
$ws.Range("I33").Value2 = 262.72726
$ws.Range("K33").Value2 = 262.72726
$ws.Range("M33").Value2 = -33.72726
$ws.Range("H40").Value2 = 3935.5454
$ws.Range("J40").Value2 = 2199.375
$ws.Range("L40").Value2 = 2199.375
$ws.Range("N40").Value2 = -2549.375
$ws.Range("H62").Value2 = 19890.584
$ws.Range("I62").Value2 = 19468.1
$ws.Range("K62").Value2 = 19468.1
$ws.Range("M62").Value2 = -18844.1
$ws.Range("H65").Value2 = 19890.584
$ws.Range("I65").Value2 = 19468.1
$ws.Range("K65").Value2 = 97340.5
$ws.Range("M65").Value2 = -94220.5
$ws.Range("H106").Value2 = 6201.5557
$ws.Range("I106").Value2 = 5759.143
$ws.Range("J106").Value2 = 7750
$ws.Range("K106").Value2 = 5759.143
$ws.Range("L106").Value2 = 7750
$ws.Range("M106").Value2 = -5128.143
$ws.Range("N106").Value2 = -9012
$ws.Range("H113").Value2 = 3129.2856
$ws.Range("I113").Value2 = 3224.75
$ws.Range("J113").Value2 = 3002
$ws.Range("K113").Value2 = 3224.75
$ws.Range("L113").Value2 = 3002
$ws.Range("M113").Value2 = 29.25
$ws.Range("N113").Value2 = -9510
$ws.Range("H137").Value2 = 38240892
$ws.Range("I137").Value2 = 90910264
$ws.Range("J137").Value2 = 2030698
$ws.Range("K137").Value2 = 272730792
$ws.Range("L137").Value2 = 6092094
$ws.Range("M137").Value2 = -272728242
$ws.Range("N137").Value2 = -6097194

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 2049.9285
$ws.Range("I45").Value2 = 2019.4
$ws.Range("K45").Value2 = 2019.4
$ws.Range("M45").Value2 = -1642.4
$ws.Range("H74").Value2 = 4033737.5
$ws.Range("I74").Value2 = 4630550.5
$ws.Range("J74").Value2 = 5250
$ws.Range("K74").Value2 = 4630550.5
$ws.Range("L74").Value2 = 5250
$ws.Range("M74").Value2 = -4629676.5
$ws.Range("N74").Value2 = -6998
$ws.Range("H77").Value2 = 4033737.5
$ws.Range("I77").Value2 = 4630550.5
$ws.Range("J77").Value2 = 5250
$ws.Range("K77").Value2 = 23152752.5
$ws.Range("L77").Value2 = 26250
$ws.Range("M77").Value2 = -23148384.5
$ws.Range("N77").Value2 = -34986

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 107425.664
$ws.Range("I31").Value2 = 209472.06
$ws.Range("J31").Value2 = 5379.2666
$ws.Range("K31").Value2 = 209472.06
$ws.Range("L31").Value2 = 5379.2666
$ws.Range("M31").Value2 = -209177.06
$ws.Range("N31").Value2 = -5969.2666
$ws.Range("H34").Value2 = 107425.664
$ws.Range("I34").Value2 = 209472.06
$ws.Range("J34").Value2 = 5379.2666
$ws.Range("K34").Value2 = 209472.06
$ws.Range("L34").Value2 = 5379.2666
$ws.Range("M34").Value2 = -209270.06
$ws.Range("N34").Value2 = -5783.2666
$ws.Range("H58").Value2 = 415201.12
$ws.Range("I58").Value2 = 618732.9
$ws.Range("J58").Value2 = 8137.7
$ws.Range("K58").Value2 = 618732.9
$ws.Range("L58").Value2 = 8137.7
$ws.Range("M58").Value2 = -618529.9
$ws.Range("N58").Value2 = -8543.700000000001
$ws.Range("H99").Value2 = 2011.4
$ws.Range("I99").Value2 = 1278.5
$ws.Range("K99").Value2 = 1278.5
$ws.Range("M99").Value2 = 219.5
$ws.Range("H107").Value2 = 1095.4286
$ws.Range("I107").Value2 = 994
$ws.Range("K107").Value2 = 994
$ws.Range("M107").Value2 = 926
$ws.Range("H126").Value2 = 2011.4
$ws.Range("I126").Value2 = 1278.5
$ws.Range("K126").Value2 = 3835.5
$ws.Range("M126").Value2 = -1365.5
$ws.Range("H132").Value2 = 75016680
$ws.Range("I132").Value2 = 111134400
$ws.Range("K132").Value2 = 333403200
$ws.Range("M132").Value2 = -333400670
$ws.Range("H134").Value2 = 23357.875
$ws.Range("I134").Value2 = 29602.75
$ws.Range("K134").Value2 = 88808.25
$ws.Range("M134").Value2 = -86273.25
$ws.Range("H136").Value2 = 415201.12
$ws.Range("I136").Value2 = 618732.9
$ws.Range("J136").Value2 = 8137.7
$ws.Range("K136").Value2 = 1856198.7
$ws.Range("L136").Value2 = 24413.1
$ws.Range("M136").Value2 = -1853648.7
$ws.Range("N136").Value2 = -29513.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value2 = 9969
$ws.Range("J101").Value2 = 9969
$ws.Range("L101").Value2 = 29907
$ws.Range("N101").Value2 = -34775
$ws.Range("H122").Value2 = 23765588
$ws.Range("J122").Value2 = 106944450
$ws.Range("L122").Value2 = 962500050
$ws.Range("N122").Value2 = -962504950

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value2 = 3270.3333
$ws.Range("I113").Value2 = 2724.4
$ws.Range("K113").Value2 = 2724.4
$ws.Range("M113").Value2 = -554.4000000000001
$ws.Range("H126").Value2 = 2780110.8
$ws.Range("I126").Value2 = 2780110.8
$ws.Range("J126").Value2 = 0
$ws.Range("K126").Value2 = 8340332.399999999
$ws.Range("L126").Value2 = 0
$ws.Range("M126").Value2 = -8337862.399999999
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 2880.8096
$ws.Range("I40").Value2 = 2880.8096
$ws.Range("J40").Value2 = 0
$ws.Range("K40").Value2 = 2880.8096
$ws.Range("L40").Value2 = 0
$ws.Range("M40").Value2 = -2744.8096
$ws.Range("H136").Value2 = 35640.08
$ws.Range("I136").Value2 = 2069.8262
$ws.Range("K136").Value2 = 6209.4786
$ws.Range("M136").Value2 = -3659.4786
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 2020.76
$ws.Range("I107").Value2 = 1198.85
$ws.Range("K107").Value2 = 3596.55
$ws.Range("M107").Value2 = -1676.55
$ws.Range("H124").Value2 = 52000
$ws.Range("J124").Value2 = 52000
$ws.Range("L124").Value2 = 52000
$ws.Range("N124").Value2 = -61820
$ws.Range("H126").Value2 = 4997.4287
$ws.Range("I126").Value2 = 4166.1665
$ws.Range("J126").Value2 = 9985
$ws.Range("K126").Value2 = 12498.4995
$ws.Range("L126").Value2 = 29955
$ws.Range("M126").Value2 = -10028.4995
$ws.Range("N126").Value2 = -34895
